$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4548.125
$ws.Range("I40").Value = 3173.6
$ws.Range("J40").Value = 5172.909
$ws.Range("K40").Value = 3173.6
$ws.Range("L40").Value = 5172.909
$ws.Range("M40").Value = -2998.6
$ws.Range("N40").Value = -5522.909
$ws.Range("H51").Value = 7799.8
$ws.Range("J51").Value = 8249.5
$ws.Range("L51").Value = 8249.5
$ws.Range("N51").Value = -9217.5
$ws.Range("H64").Value = 3948.75
$ws.Range("I64").Value = 3795
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 3795
$ws.Range("L64").Value = 4000
$ws.Range("M64").Value = -3547
$ws.Range("N64").Value = -4496
$ws.Range("H67").Value = 3948.75
$ws.Range("I67").Value = 3795
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 3795
$ws.Range("L67").Value = 4000
$ws.Range("M67").Value = -2937
$ws.Range("N67").Value = -5716
$ws.Range("H74").Value = 7000
$ws.Range("I74").Value = 7000
$ws.Range("K74").Value = 7000
$ws.Range("M74").Value = -6064
$ws.Range("H77").Value = 7000
$ws.Range("I77").Value = 7000
$ws.Range("K77").Value = 35000
$ws.Range("M77").Value = -30320
$ws.Range("H116").Value = 4404.3335
$ws.Range("I116").Value = 4185.4
$ws.Range("K116").Value = 4185.4
$ws.Range("M116").Value = -743.3999999999996

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 56198.668
$ws.Range("J24").Value = 56198.668
$ws.Range("L24").Value = 56198.668
$ws.Range("N24").Value = -56946.668
$ws.Range("H94").Value = 50880.5
$ws.Range("J94").Value = 50880.5
$ws.Range("L94").Value = 50880.5
$ws.Range("N94").Value = -52682.5
$ws.Range("H97").Value = 2270.1667
$ws.Range("I97").Value = 2057.8333
$ws.Range("K97").Value = 2057.8333
$ws.Range("M97").Value = -1561.8333
$ws.Range("H100").Value = 56198.668
$ws.Range("J100").Value = 56198.668
$ws.Range("L100").Value = 56198.668
$ws.Range("N100").Value = -58362.668
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H106").Value = 50370
$ws.Range("J106").Value = 50370
$ws.Range("L106").Value = 50370
$ws.Range("N106").Value = -52894
$ws.Range("H122").Value = 1408.4
$ws.Range("I122").Value = 1408.4
$ws.Range("K122").Value = 4225.200000000001
$ws.Range("M122").Value = -1775.200000000001

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 870.53845
$ws.Range("I20").Value = 793.1111
$ws.Range("K20").Value = 793.1111
$ws.Range("M20").Value = -546.1111
$ws.Range("H86").Value = 8622
$ws.Range("I86").Value = 9994
$ws.Range("J86").Value = 8164.6665
$ws.Range("K86").Value = 9994
$ws.Range("L86").Value = 8164.6665
$ws.Range("M86").Value = -8871
$ws.Range("N86").Value = -10410.6665
$ws.Range("H89").Value = 8622
$ws.Range("I89").Value = 9994
$ws.Range("J89").Value = 8164.6665
$ws.Range("K89").Value = 49970
$ws.Range("L89").Value = 40823.3325
$ws.Range("M89").Value = -44354
$ws.Range("N89").Value = -52055.3325
$ws.Range("H99").Value = 2435.3635
$ws.Range("I99").Value = 1975
$ws.Range("J99").Value = 3663
$ws.Range("K99").Value = 1975
$ws.Range("L99").Value = 3663
$ws.Range("M99").Value = -477
$ws.Range("N99").Value = -6659
$ws.Range("H105").Value = 1060.5
$ws.Range("I105").Value = 1060.5
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1060.5
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 686.5
$ws.Range("N105").ClearContents()

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 428.5
$ws.Range("I5").Value = 235.33333
$ws.Range("K5").Value = 235.33333
$ws.Range("M5").Value = -123.33333
$ws.Range("H31").Value = 1499.3334
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H34").Value = 1499.3334
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H99").Value = 2558.5386
$ws.Range("I99").Value = 2679.1667
$ws.Range("J99").Value = 1111
$ws.Range("K99").Value = 2679.1667
$ws.Range("L99").Value = 1111
$ws.Range("M99").Value = -1181.1667
$ws.Range("N99").Value = -4107
$ws.Range("H107").Value = 1171.5
$ws.Range("I107").Value = 865.64
$ws.Range("K107").Value = 865.64
$ws.Range("M107").Value = 1054.36
$ws.Range("H122").Value = 955.4
$ws.Range("J122").Value = 250
$ws.Range("L122").Value = 750
$ws.Range("N122").Value = -5650
$ws.Range("H126").Value = 2558.5386
$ws.Range("I126").Value = 2679.1667
$ws.Range("J126").Value = 1111
$ws.Range("K126").Value = 8037.500100000001
$ws.Range("L126").Value = 3333
$ws.Range("M126").Value = -5567.500100000001
$ws.Range("N126").Value = -8273
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1842
$ws.Range("I18").Value = 1447
$ws.Range("K18").Value = 4341
$ws.Range("M18").Value = -4172
$ws.Range("H22").Value = 65675.5
$ws.Range("I22").Value = 84233.336
$ws.Range("K22").Value = 252700.008
$ws.Range("M22").Value = -252531.008
$ws.Range("H27").Value = 65675.5
$ws.Range("I27").Value = 84233.336
$ws.Range("K27").Value = 252700.008
$ws.Range("M27").Value = -252598.008

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 12101.6
$ws.Range("I43").Value = 11504
$ws.Range("J43").Value = 12500
$ws.Range("K43").Value = 11504
$ws.Range("L43").Value = 12500
$ws.Range("M43").Value = -11353
$ws.Range("N43").Value = -12802
$ws.Range("H80").Value = 7501.5
$ws.Range("J80").Value = 8835.333000000001
$ws.Range("L80").Value = 8835.333000000001
$ws.Range("N80").Value = -10831.333
$ws.Range("H83").Value = 7501.5
$ws.Range("J83").Value = 8835.333000000001
$ws.Range("L83").Value = 44176.665
$ws.Range("N83").Value = -54160.665
$ws.Range("H97").Value = 4068.1428
$ws.Range("I97").Value = 3747
$ws.Range("J97").Value = 4496.3335
$ws.Range("K97").Value = 3747
$ws.Range("L97").Value = 4496.3335
$ws.Range("M97").Value = -3251
$ws.Range("N97").Value = -5488.3335

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H100").Value = 12345
$ws.Range("I100").Value = 12345
$ws.Range("K100").Value = 12345
$ws.Range("M100").Value = -11804
$ws.Range("H104").Value = 47370
$ws.Range("J104").Value = 47370
$ws.Range("L104").Value = 47370
$ws.Range("N104").Value = -54358

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H122").Value = 3937.5
$ws.Range("I122").Value = 4000
$ws.Range("K122").Value = 12000
$ws.Range("M122").Value = -9550

